$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 66)
$ws.Range("A66").Value = "NEKO-MIMI"
$ws.Range("B66").Value = 1993
$ws.Range("D66").Value = "https://letterboxd.com/film/neko-mimi/"
$ws.Range("E66").Value = "Missing or blank fields: Genres"

# Set the row height for the newly added row
$ws.Range("A66:E66").RowHeight = 16

# Update the selected cell / view state
$ws.Range("O57").Select()
